$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update E2: 1970.64 -> 2000
$ws.Range("E2").Value = 2000

# 2. Update D13/E13 (same date/currency/price row)
$ws.Range("D13").Value = 43834.38366132
$ws.Range("E13").Value = 2563.87

# 3. Insert a new row at 14, shifting existing rows 14-45 down to 15-46
$ws.Rows("14:14").Insert()

# 4. Populate the newly inserted row 14 with its data
$ws.Range("A14").Value = 45146.03059027778
$ws.Range("B14").Value = "HBAR"
$ws.Range("C14").Value = 0.05909090908667258
$ws.Range("D14").Value = 53.64615385
$ws.Range("E14").Value = 3.17

# 5. Fix D38/E38 (originally row 37 pre-insert, now shifted to row 38) to the new values
$ws.Range("D38").Value = 6661.40079847
$ws.Range("E38").Value = 506.82
